$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: paragraph "[module 06: styling CSS components using css
# module]" -> "[module 06: styling CSS components] using CSS module]"
# split across five runs.
# ---------------------------------------------------------------------
$full = $d.Content.Text
$needle1 = "styling CSS components using css module]"
$idx1 = $full.IndexOf($needle1)
$rng1 = $d.Range($idx1, $idx1 + $needle1.Length)

$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
'<w:body><w:p>' + `
'<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>styling CSS components]</w:t></w:r>' + `
'<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
'<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">using </w:t></w:r>' + `
'<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>CSS</w:t></w:r>' + `
'<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> module</w:t></w:r>' + `
'</w:p></w:body></w:document>' + `
'</pkg:xmlData></pkg:part></pkg:package>'

$rng1.InsertXML($xml1)

# ---------------------------------------------------------------------
# Change 2: paragraph "Let's see the button class.<br>css module gives
# a unique value to the button, <br>component name => our class name
# => unique class name" gets "css" wrapped in a spell-check proofErr
# pair and split into its own run.
# ---------------------------------------------------------------------
$p7 = $d.Paragraphs(7)
$p7Range = $p7.Range
$p7RangeNoMark = $d.Range($p7Range.Start, $p7Range.End - 1)

$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
'<w:body><w:p>' + `
'<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Let</w:t></w:r>' + `
'<w:r w:rsidR="00FC7173"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>’</w:t></w:r>' + `
'<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>s see the button class.</w:t></w:r>' + `
'<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br/></w:r>' + `
'<w:proofErr w:type="spellStart"/>' + `
'<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>css</w:t></w:r>' + `
'<w:proofErr w:type="spellEnd"/>' + `
'<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> module gives a unique value to the button, </w:t></w:r>' + `
'<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br/><w:t xml:space="preserve">component name =&gt; </w:t></w:r>' + `
'<w:r w:rsidR="00BD789F"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>our class name</w:t></w:r>' + `
'<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> =&gt; unique class name</w:t></w:r>' + `
'</w:p></w:body></w:document>' + `
'</pkg:xmlData></pkg:part></pkg:package>'

$p7RangeNoMark.InsertXML($xml2)
